$wb = $excel.ActiveWorkbook

# Rename the "Include from ..." sheets to "Include #N"
$wb.Worksheets.Item("Include from Tipo de Observaç").Name = "Include #0"
$wb.Worksheets.Item("Include from Tabela de proced").Name = "Include #1"
$wb.Worksheets.Item("Include from LOINC").Name = "Include #2"

# Update the Contact value on the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B10").Value = "null (http://www.saude.gov.br)"
